$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 232, shifting rows 232:348 down to 233:349
$ws.Rows.Item(232).Insert()

# Populate the newly inserted row 232 with the new data record
$ws.Range("A232").Value = 5
$ws.Range("B232").Value = "Macroferia Regional de Talca"
$ws.Range("C232").Value = "Maule"
$ws.Range("D232").Value = 44917
$ws.Range("D232").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E232").Value = 7
$ws.Range("F232").Value = 100112009
$ws.Range("G232").Value = "Acelga"
$ws.Range("H232").Value = "Sin especificar"
$ws.Range("I232").Value = "Primera"
$ws.Range("J232").Value = 250
$ws.Range("K232").Value = 3000
$ws.Range("L232").Value = 3000
$ws.Range("M232").Value = 3000
$ws.Range("N232").Value = "$/docena de atados (4 kilos)"
$ws.Range("O232").Value = "Región del Maule"
$ws.Range("P232").Value = 750
$ws.Range("Q232").Value = 4
$ws.Range("R232").Value = "Hortaliza"
